$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the price column keeps its textual representation (avoid Excel
# auto-converting numeric-looking strings like "330.55" into numbers).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '27.617.32'
$ws.Range("E2").Value = '  -1.83%  '
$ws.Range("D3").Value = '1.737.91'
$ws.Range("E3").Value = '  -2.56%  '
$ws.Range("D4").Value = '1.008'
$ws.Range("E4").Value = '  +0.87%  '
$ws.Range("D5").Value = '330.55'
$ws.Range("E6").Value = '  +0.62%  '
$ws.Range("D7").Value = '0.3861'
$ws.Range("E7").Value = '  +0.99%  '
$ws.Range("D8").Value = '0.3339'
$ws.Range("E8").Value = '  -2.64%  '
$ws.Range("D9").Value = '45.29'
$ws.Range("E9").Value = '  -5.08%  '
$ws.Range("D10").Value = '1.091'
$ws.Range("E10").Value = '  -5.12%  '
$ws.Range("D11").Value = '0.07095'
$ws.Range("E11").Value = '  -4.65%  '
$ws.Range("E12").Value = '  +0.78%  '
$ws.Range("D13").Value = '21.70'
$ws.Range("E13").Value = '  -6.22%  '
$ws.Range("D14").Value = '6.047'
$ws.Range("E14").Value = '  -5.56%  '
$ws.Range("D15").Value = '1.742.18'
$ws.Range("E15").Value = '  -2.44%  '
$ws.Range("D16").Value = '6.892'
$ws.Range("E16").Value = '  -3.40%  '
$ws.Range("D17").Value = '0.00001038'
$ws.Range("E17").Value = '  -4.17%  '
$ws.Range("D18").Value = '0.06571'
$ws.Range("E18").Value = '  -1.14%  '
$ws.Range("E19").Value = '  +0.67%  '
$ws.Range("D20").Value = '78.11'
$ws.Range("E20").Value = '  -5.70%  '
$ws.Range("D21").Value = '16.47'
$ws.Range("E21").Value = '  -5.79%  '
$ws.Range("D22").Value = '6.124'
$ws.Range("E22").Value = '  -4.89%  '
$ws.Range("D23").Value = '27.646.30'
$ws.Range("E23").Value = '  -1.69%  '
$ws.Range("D24").Value = '11.37'
$ws.Range("E24").Value = '  -5.81%  '
$ws.Range("E25").Value = '  +0.37%  '
$ws.Range("D26").Value = '153.09'
$ws.Range("E26").Value = '  -0.66%  '
$ws.Range("D27").Value = '19.45'
$ws.Range("E27").Value = '  -7.01%  '
$ws.Range("D28").Value = '2.245'
$ws.Range("E28").Value = '  -7.94%  '
$ws.Range("D29").Value = '1.939.25'
$ws.Range("E29").Value = '  -2.54%  '
$ws.Range("D30").Value = '1.240'
$ws.Range("E30").Value = '  -14.31%  '
$ws.Range("D31").Value = '126.63'
$ws.Range("E31").Value = '  -5.91%  '
$ws.Range("D32").Value = '4.029'
$ws.Range("E32").Value = '  +1.87%  '
$ws.Range("D33").Value = '5.675'
$ws.Range("E33").Value = '  -7.87%  '
$ws.Range("D34").Value = '0.08623'
$ws.Range("E34").Value = '  -1.08%  '
$ws.Range("D35").Value = '11.79'
$ws.Range("E35").Value = '  -7.94%  '
$ws.Range("B36").Value = 'WEMIXTOKEN'
$ws.Range("C36").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D36").Value = '1.497'
$ws.Range("E36").Value = '  -1.32%  '
$ws.Range("B37").Value = 'InternetComputer(DFINITY)'
$ws.Range("C37").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D37").Value = '5.041'
$ws.Range("E37").Value = '  -5.45%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = '0.02222'
$ws.Range("E38").Value = '  -8.47%  '
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").Value = '0.05975'
$ws.Range("E39").Value = '  -5.66%  '
$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D40").Value = '0.6309'
$ws.Range("E40").Value = '  -8.05%  '
$ws.Range("D41").Value = '0.2063'
$ws.Range("E41").Value = '  -5.90%  '
$ws.Range("D42").Value = '1.183'
$ws.Range("E42").Value = '  -4.86%  '
$ws.Range("E43").Value = '  +0.68%  '
$ws.Range("D44").Value = '7.816'
$ws.Range("E44").Value = '  -6.44%  '
$ws.Range("E45").Value = '  -4.82%  '
$ws.Range("E46").Value = '  -1.40%  '
$ws.Range("D47").Value = '0.5850'
$ws.Range("E47").Value = '  -7.39%  '
$ws.Range("D48").Value = '124.38'
$ws.Range("E48").Value = '  -6.00%  '
$ws.Range("D49").Value = '1.942'
$ws.Range("E49").Value = '  -7.32%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").Value = '0.06905'
$ws.Range("E50").Value = '  -7.27%  '
$ws.Range("B51").Value = 'EOS'
$ws.Range("C51").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D51").Value = '1.134'
$ws.Range("E51").Value = '  -2.00%  '
